$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the defect/fault descriptions in column C for rows 2-29
$ws.Range("C2").Value  = "Механические повреждения, Принтер"
$ws.Range("C4").Value  = "Залит, Принтер"
$ws.Range("C5").Value  = "Следы жизнедеятельности, GPRS"
$ws.Range("C6").Value  = "Неустранимые загрязнения, Принтер"
$ws.Range("C7").Value  = "Следы жизнедеятельности, GPRS"
$ws.Range("C8").Value  = "Залит, Tamper"
$ws.Range("C9").Value  = "Залит, Дефект клавиатуры"
$ws.Range("C10").Value = "Неустранимые загрязнения, Дефект клавиатуры"
$ws.Range("C11").Value = "Механические повреждения, Принтер"
$ws.Range("C12").Value = "Механические повреждения, CTLS"
$ws.Range("C14").Value = "Залит, Принтер"
$ws.Range("C15").Value = "Неустранимые загрязнения, Дефект клавиатуры"
$ws.Range("C16").Value = "Залит, Принтер"
$ws.Range("C17").Value = "Залит, Порт питания"
$ws.Range("C18").Value = "Механические повреждения, Tamper"
$ws.Range("C19").Value = "Механические повреждения, GPRS"
$ws.Range("C20").Value = "Механические повреждения, Дефект экрана"
$ws.Range("C21").Value = "Неустранимые загрязнения, Принтер"
$ws.Range("C23").Value = "Залит, CTLS"
$ws.Range("C24").Value = "Залит, Alert"
$ws.Range("C25").Value = "Механические повреждения, CTLS"
$ws.Range("C26").Value = "Залит, CTLS"
$ws.Range("C27").Value = "Неустранимые загрязнения, Дефект клавиатуры"
$ws.Range("C28").Value = "Залит, GPRS"
$ws.Range("C29").Value = "Механические повреждения, CTLS"

# Add a new row of data at row 30
$ws.Range("A30").Value = "Castles"
$ws.Range("B30").Value = "12345345687u"
$ws.Range("C30").Value = "Неустранимые загрязнения"

# Match the author's style formatting used by other Model/S-N cells
$ws.Range("A30").Style = "Обычный"
$ws.Range("B30").Style = "Обычный"

$ws.Range("A30:B30").NumberFormat = "@"

# Update selection to reflect the new extended data range
$ws.Range("C2:C34").Select()
